# Append/update the latest ランサーズ job-listing row (row 2) and widen a
# couple of columns to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column widths -------------------------------------------------------
# Excel's ColumnWidth setter pads by ~0.8333 chars (the default-font
# padding) before it lands in the saved `width` attribute, so back that
# off here to land on exactly 34 / 27 in the saved file.
$ws.Columns.Item(2).ColumnWidth = 33.166666666666664   # B: 20 -> 34
$ws.Columns.Item(4).ColumnWidth = 26.166666666666668   # D: 26 -> 27

# --- Row 2 cell values -----------------------------------------------------
$ws.Range("A2").Value = "2025-12-14 18:24:59"
$ws.Range("B2").Value = "Base無在庫ツール作成 経験者のみ募集 実績提示をお願いします"
$ws.Range("D2").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("G2").Value = 73
$ws.Range("H2").Value = "◆ツール"

# --- F2: new URL text + matching hyperlink ---------------------------------
$ws.Range("F2").Hyperlinks.Delete()
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5453611"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5453611")
